$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("MySingleFunction")
$newSheet = $ws3.Copy([Type]::Missing, $ws3)
$wb.Worksheets.Item("MySingleFunction (2)").Name = "Alt1"
